# The slide master and every slide layout carry a cached copy of the
# auto-updating "date" footer field (type="datetimeFigureOut"). The
# presentation was re-saved a day later (16/3/2025 -> 17/3/2025), so
# PowerPoint refreshed the cached display text for that field on the
# master and on every layout. Reproduce that by updating the "Date
# Placeholder" shape's text wherever it appears.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "16/3/2025") {
                $tr.Text = $newText
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "17/3/2025"

# Every slide layout under the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "17/3/2025"
}
